$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (and would lose formatting, e.g. "0.360" -> 0.36, or trailing-zero loss).
# For those we force a Text number format, assign the value, then restore the
# default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = "57.541.59"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "2.327.61"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.355.65"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.31%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.774.43"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "57.689.19"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D18").Value = "2.341.89"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +6.05%  "
$ws.Range("E29").Value = "  +5.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +16.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +7.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "149.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "289.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.05%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0505"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.562"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.380"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
